$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.576.12"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "2.320.05"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'496.46"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "'128.92"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  -3.86%  "
$ws.Range("D9").Value = "2.312.23"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'0.0941"
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("D13").Value = "'0.316"
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").Value = "2.735.76"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "'21.34"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "54.563.60"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "2.329.34"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").Value = "'9.69"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").Value = "'3.97"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "'304.14"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'64.31"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'0.367"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'0.143"
$ws.Range("E27").Value = "  -4.87%  "
$ws.Range("D28").Value = "'7.11"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").Value = "'167.51"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").Value = "0.0₃0699"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D33").Value = "'5.75"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("D36").Value = "'17.56"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("D38").Value = "'0.849"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("D40").Value = "'35.79"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -3.94%  "
$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "'124.40"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("D45").Value = "'4.70"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0885"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.546"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").Value = "'238.06"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "'0.0474"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").Value = "'0.0204"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "'16.49"
$ws.Range("E51").Value = "  -2.55%  "
